$wb = $excel.ActiveWorkbook

# Source sheets to copy existing formatting / shared-string values from,
# so new cells reuse the same styles & shared strings as the rest of the
# workbook instead of Excel inventing new ones (e.g. auto-detecting dates).
$debtSheet = $wb.Worksheets.Item(7)     # "債務" - same row layout as the new sheet
$landSheet = $wb.Worksheets.Item(1)     # "土地" - holds the literal "2012-02-01" text in K2

# Add the new worksheet as the last tab in the workbook.
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "事業投資"

# Bring over the header-row style (bold/centered/bordered) and the
# data-row style used by every other sheet in this workbook.
$debtSheet.Range("B1:N1").Copy()
$newSheet.Range("B1:N1").PasteSpecial(-4122)  # xlPasteFormats

$debtSheet.Range("A2:N2").Copy()
$newSheet.Range("A2:N2").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "owner"
$newSheet.Range("C1").Value = "company"
$newSheet.Range("D1").Value = "address"
$newSheet.Range("E1").Value = "total"
$newSheet.Range("F1").Value = "register_date"
$newSheet.Range("G1").Value = "register_reason"
$newSheet.Range("H1").Value = "property_category"
$newSheet.Range("I1").Value = "category"
$newSheet.Range("J1").Value = "date"
$newSheet.Range("K1").Value = "legislator_name"
$newSheet.Range("L1").Value = "legislator_id"
$newSheet.Range("M1").Value = "source_file"
$newSheet.Range("N1").Value = "index"

# Data row
$newSheet.Range("A2").Value = 109
$newSheet.Range("B2").Value = "吳宜臻"
$newSheet.Range("C2").Value = "志遠法律事務所"
$newSheet.Range("D2").Value = "臺北市大安區仁愛路三段98號4樓"
$newSheet.Range("E2").Value = 1500000
$newSheet.Range("F2").Value = "100年05月30日"
$newSheet.Range("G2").Value = "個人投資"
$newSheet.Range("H2").Value = "investment"
$newSheet.Range("I2").Value = "normal"

# "2012-02-01" must stay a literal text value (not get auto-converted to a
# date serial). Copy it from a cell elsewhere in the workbook that already
# holds this exact text ("土地"!K2) so it lands as plain text/shared string.
$landSheet.Range("K2").Copy()
$newSheet.Range("J2").PasteSpecial(-4163)  # xlPasteValues

$newSheet.Range("K2").Value = "吳宜臻"
$newSheet.Range("L2").Value = 1735
$newSheet.Range("M2").Value = "tmp2691"
$newSheet.Range("N2").Value = 109
